# Add new localization entries for "Recent Actions Filter" feature
# to the "General" worksheet of the Emmersive localization workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# New row 120: em_ui_filter / 最近の会話フィルター / Recent Actions Filter
# New row 121: em_ui_add / 追加  / Add
# (values are set in this specific order so that new shared-string
#  table entries are appended in the same sequence as the source edit)
$ws.Range("A120").Value = "em_ui_filter"
$ws.Range("A121").Value = "em_ui_add"
$ws.Range("D121").Value = "Add"
$ws.Range("C121").Value = "追加 "
$ws.Range("C120").Value = "最近の会話フィルター"
$ws.Range("D120").Value = "Recent Actions Filter"

# Update selection/view to match the authored state
$ws.Range("D123").Select()
